$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4452382786404918
$ws.Range("D2").Value = 14025372236.59455
$ws.Range("G2").Value = 0.3
$ws.Range("I2").Value = 7587630754.078597
$ws.Range("M2").Value = 11159709000

# Row 3
$ws.Range("B3").Value = 0.4452382786404918
$ws.Range("D3").Value = 14025372236.59455
$ws.Range("G3").Value = 0.3
$ws.Range("I3").Value = 7587630754.078597
$ws.Range("M3").Value = 11159709000
